$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "67.418.08"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -1.82%  "
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.680.52"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.07%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "598.40"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "163.72"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("E7").Value = "  -0.02%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.545"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.679.77"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("E10").Value = "  +1.77%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.358"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.26%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "5.22"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -1.37%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.82"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.89%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.170.00"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.64%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0000184"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.96%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.439.36"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -1.75%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.657.35"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -1.52%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.64"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "362.92"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -1.37%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.50"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.65%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.40"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "4.81"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.51%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.03"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -4.13%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "71.60"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -5.05%  "
$ws.Range("E26").Value = "  +0.10%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "10.00"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.52%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.817.63"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("E29").Value = "  -1.79%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "550.27"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -6.26%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.97"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -3.84%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.38"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -3.45%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.98%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "0.130"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.12%  "
$ws.Range("E36").Value = "  -0.01%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "1.58"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -4.79%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "19.56"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -2.60%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "156.59"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -2.79%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.372"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.17%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.84"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -3.30%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.27"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -2.61%  "
$ws.Range("E43").Value = "  +0.28%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.53"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("E45").Value = "  +0.04%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "40.36"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.40%  "
$ws.Range("E47").Value = "  -6.15%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "0.586"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -4.37%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "152.74"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.60%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "3.84"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.85%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.73"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -3.79%  "
